$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-11 Sunday" "2026-01-12 Monday"

Replace-Text "61÷3=20, 1" "82÷7=11, 5"
Replace-Text "39÷7=5, 4" "51÷2=25, 1"
Replace-Text "32÷3=10, 2" "16÷4=4, 0"
Replace-Text "55÷8=6, 7" "49÷9=5, 4"
Replace-Text "73÷2=36, 1" "86÷8=10, 6"

Replace-Text "29÷3=9, 2" "32÷7=4, 4"
Replace-Text "40÷3=13, 1" "25÷4=6, 1"
Replace-Text "96÷9=10, 6" "63÷7=9, 0"
Replace-Text "54÷6=9, 0" "74÷6=12, 2"
Replace-Text "37÷8=4, 5" "61÷2=30, 1"

Replace-Text "93÷8=11, 5" "61÷4=15, 1"
Replace-Text "76÷2=38, 0" "98÷8=12, 2"
Replace-Text "28÷2=14, 0" "56÷5=11, 1"
Replace-Text "94÷6=15, 4" "48÷8=6, 0"
Replace-Text "86÷4=21, 2" "33÷7=4, 5"

Replace-Text "15÷7=2, 1" "48÷5=9, 3"
Replace-Text "38÷8=4, 6" "66÷2=33, 0"
Replace-Text "43÷7=6, 1" "86÷6=14, 2"
Replace-Text "61÷6=10, 1" "31÷3=10, 1"
Replace-Text "29÷2=14, 1" "92÷5=18, 2"

Replace-Text "32÷4=8, 0" "67÷8=8, 3"
Replace-Text "67÷2=33, 1" "38÷7=5, 3"
Replace-Text "26÷8=3, 2" "79÷8=9, 7"
Replace-Text "59÷3=19, 2" "74÷4=18, 2"
Replace-Text "78÷8=9, 6" "99÷8=12, 3"
